# Updated symbol list on Tue Dec 27 13:31:05 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # The "Price" column stores plain text that happens to look numeric
    # (e.g. "243.71"). Force Excel to keep it as text (quote-prefix trick)
    # and then reset the cell style back to Normal so no stray
    # number-format style gets attached to the cell.
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Price tweaks that don't involve row movement
Set-TextValue $ws.Range("D2") "243.71"
Set-TextValue $ws.Range("D3") "24.10"
Set-TextValue $ws.Range("D4") "5.374"
Set-TextValue $ws.Range("D6") "3.412"
Set-TextValue $ws.Range("D8") "0.8105"
Set-TextValue $ws.Range("D9") "0.9488"

# Rows 10-18: "One" moved to the top (row 10), shifting WazirX..CoinExToken down by one row
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D10") "0.0005959"
$ws.Range("E10").Value = "9OneONE"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D11") "0.1422"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.07427"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D13") "0.03107"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D14") "0.03041"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D15") "0.09342"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D16") "3.865"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D17") "0.001576"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D18") "0.04711"
$ws.Range("E18").Value = "17CoinExTokenCET"

# Remaining isolated price / label tweaks
Set-TextValue $ws.Range("D19") "0.005947"
Set-TextValue $ws.Range("D20") "0.001244"
Set-TextValue $ws.Range("D22") "0.00007998"
$ws.Range("E22").Value = "21NitroExNTXWorstin24h"
Set-TextValue $ws.Range("D25") "0.3223"
Set-TextValue $ws.Range("D27") "0.0002652"
Set-TextValue $ws.Range("D40") "0.03904"
Set-TextValue $ws.Range("D41") "0.006369"
Set-TextValue $ws.Range("D42") "0.1071"
Set-TextValue $ws.Range("D43") "0.003199"
Set-TextValue $ws.Range("D44") "0.008497"
Set-TextValue $ws.Range("D45") "0.00005208"
Set-TextValue $ws.Range("D47") "0.7198"
